# Add the newly-discovered article URL to the "遠隔会議" (remote
# conferencing) sheet. Every other sheet in this robot-news scraper
# workbook is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("遠隔会議")

$ws.Range("A3").Value = "https://robotstart.info/2018/01/09/fairydevices-tumbler.html"
# New rows pick up the column's default style (style id 1) on write;
# the source row has no explicit style, so reset it to Normal/General.
$ws.Range("A3").Style = "Normal"
